# Apply the "Bracket" sheet pick updates.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bracket")

$ws.Range("E4").Value  = "Giant striped mongoose"
$ws.Range("F8").Value  = "Chequered elephant shrew"
$ws.Range("D10").Value = "Side-striped jackal"
$ws.Range("N10").Value = "Trapdoor Spider"
$ws.Range("M12").Value = "Lungfish"
$ws.Range("N14").Value = "Lungfish"
$ws.Range("G16").Value = "Kudu"
$ws.Range("K16").Value = "Golden Eagle"
$ws.Range("N18").Value = "Goanna"
$ws.Range("E20").Value = "Striped hyena"
$ws.Range("M20").Value = "Bee"
$ws.Range("N22").Value = "Bee"
$ws.Range("F24").Value = "Kudu"
$ws.Range("L24").Value = "Dung Beetle"
$ws.Range("D26").Value = "Numbat"
$ws.Range("E28").Value = "Kudu"
$ws.Range("M28").Value = "Dung Beetle"
$ws.Range("D30").Value = "Kudu"
$ws.Range("N30").Value = "Dung Beetle"
$ws.Range("I32").Value = "Greater Rhea"
$ws.Range("J32").Value = "Greater Rhea"
$ws.Range("D34").Value = "Sea Otter"
$ws.Range("M36").Value = "Caspian Terns"
$ws.Range("L40").Value = "Siamang"
$ws.Range("D42").Value = "Sibree Dwarf Lemur"
$ws.Range("N42").Value = "Pacific Spiny Lumpsucker"
$ws.Range("K48").Value = "Greater Rhea"
$ws.Range("E52").Value = "Dik Dik"
$ws.Range("M52").Value = "Spotted sandpiper"
$ws.Range("D54").Value = "Dik Dik"
$ws.Range("N54").Value = "Giant Water Bug"
$ws.Range("D58").Value = "Thor's Hero Shrew"
$ws.Range("N58").Value = "Greater Flamingo"
